$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose "curated dimension" became a "curated measure":
# temporalidad (A), cif (B), mes-nombre (I), dias-duracion-contrato (K), sexo (N)
$cols = @("A", "B", "I", "K", "N")

foreach ($col in $cols) {
    # Row 2: iaest-dimension:<name> -> iaest-measure:<name>
    $cell2 = $ws.Range($col + "2")
    $cell2.Value = $cell2.Value2 -replace "^iaest-dimension:", "iaest-measure:"

    # Row 3: dim -> medida
    $ws.Range($col + "3").Value = "medida"

    # Row 4: skos:Concept -> xsd:int
    $ws.Range($col + "4").Value = "xsd:int"

    # Row 5: mapping-*.xlsx reference removed (cell removed entirely)
    $ws.Range($col + "5").Clear()
}
